$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - Vitor Ito, phone 11966548087, payment (idPagamento) still missing
$ws.Range("A14").Value = "Vitor Ito"
$ws.Range("B14").Value = 1578424633
$ws.Range("C14").Value = "'11966548087"
$ws.Range("D14").Value = "'"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 7
$ws.Range("L14").Value = 8
$ws.Range("M14").Value = 9
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = "Não"

# Row 15 - Vitor Ito, phone 119988776655, payment (idPagamento) still missing
$ws.Range("A15").Value = "Vitor Ito"
$ws.Range("B15").Value = 1578424633
$ws.Range("C15").Value = "'119988776655"
$ws.Range("D15").Value = "'"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 7
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 16
$ws.Range("L15").Value = 17
$ws.Range("M15").Value = 21
$ws.Range("N15").Value = 22
$ws.Range("O15").Value = "Não"
